# Append a new row (58) of logged data to each of the four worksheets,
# mirroring the existing row layout:
#   A: timestamp (date/time, same style as the cell above)
#   B: inlineStr  "总长"      hex bytes
#   C: inlineStr  "ID"        hex bytes
#   D: inlineStr  "实际长度"  hex bytes
#   E: inlineStr  "和校验"    hex bytes
#   F: 总长_DEC      (n)
#   G: ID_DEC        (n)
#   H: 实际长度_DEC  (n)
#   I: 和校验_DEC    (n)

$wb = $excel.ActiveWorkbook

$newRow = 58
$srcRow = 57

$rowsData = @(
    @{ Sheet = 1; A = [double]"45753.3650134375";  B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x6e"; E = "0xd"; F = 400; G = [double]"5.68631262647114e+23"; H = 366; I = 13 },
    @{ Sheet = 2; A = [double]"45753.21754209491"; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x6e"; E = "0xe"; F = 400; G = [double]"5.68631262647114e+23"; H = 366; I = 14 },
    @{ Sheet = 3; A = [double]"45753.35777591435"; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x6e"; E = "0x3"; F = 400; G = [double]"5.68631262647114e+23"; H = 366; I = 3 },
    @{ Sheet = 4; A = [double]"45753.4175328588";  B = "0x01,0x90"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x01,0x6e"; E = "0x3"; F = 400; G = [double]"9.85046333984776e+23"; H = 366; I = 3 }
)

foreach ($rd in $rowsData) {
    $ws = $wb.Worksheets.Item($rd.Sheet)

    $ws.Cells.Item($newRow, 1).Value = $rd.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($srcRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $rd.B
    $ws.Cells.Item($newRow, 3).Value = $rd.C
    $ws.Cells.Item($newRow, 4).Value = $rd.D
    $ws.Cells.Item($newRow, 5).Value = $rd.E

    $ws.Cells.Item($newRow, 6).Value = $rd.F
    $ws.Cells.Item($newRow, 7).Value = $rd.G
    $ws.Cells.Item($newRow, 8).Value = $rd.H
    $ws.Cells.Item($newRow, 9).Value = $rd.I
}
